$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking price values must stay as exact text (preserve trailing zeros / formatting),
# so force the Text number format before assigning these.
$textCells = @('D5', 'D6', 'D7', 'D9', 'D10', 'D11', 'D12', 'D13', 'D15', 'D16', 'D20', 'D21', 'D23', 'D25', 'D26', 'D27', 'D30', 'D31', 'D32', 'D33', 'D34', 'D36', 'D37', 'D39', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '42.153.35'
$ws.Range("E2").Value = '  -1.80%  '
$ws.Range("D3").Value = '2.243.60'
$ws.Range("E3").Value = '  -2.07%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '247.16'
$ws.Range("E5").Value = '  -2.27%  '
$ws.Range("D6").Value = '0.623'
$ws.Range("E6").Value = '  -3.86%  '
$ws.Range("D7").Value = '74.08'
$ws.Range("E7").Value = '  -1.29%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").Value = '0.617'
$ws.Range("E9").Value = '  -4.84%  '
$ws.Range("D10").Value = '40.90'
$ws.Range("E10").Value = '  +4.70%  '
$ws.Range("D11").Value = '0.0940'
$ws.Range("E11").Value = '  -4.58%  '
$ws.Range("D12").Value = '7.08'
$ws.Range("E12").Value = '  -6.19%  '
$ws.Range("D13").Value = '0.102'
$ws.Range("E13").Value = '  -3.39%  '
$ws.Range("D14").Value = '2.581.50'
$ws.Range("E14").Value = '  -1.96%  '
$ws.Range("D15").Value = '14.44'
$ws.Range("E15").Value = '  -4.62%  '
$ws.Range("D16").Value = '0.852'
$ws.Range("E16").Value = '  -2.48%  '
$ws.Range("D17").Value = '2.252.22'
$ws.Range("E17").Value = '  -1.91%  '
$ws.Range("D18").Value = '42.034.10'
$ws.Range("E18").Value = '  -1.79%  '
$ws.Range("E19").Value = '  -2.23%  '
$ws.Range("D20").Value = '6.11'
$ws.Range("E20").Value = '  -2.19%  '
$ws.Range("D21").Value = '71.74'
$ws.Range("E21").Value = '  -0.83%  '
$ws.Range("E22").Value = '  +6.60%  '
$ws.Range("D23").Value = '229.50'
$ws.Range("E23").Value = '  -3.36%  '
$ws.Range("D25").Value = '11.08'
$ws.Range("E25").Value = '  -2.38%  '
$ws.Range("D26").Value = '3.56'
$ws.Range("E26").Value = '  -8.24%  '
$ws.Range("D27").Value = '7.63'
$ws.Range("E27").Value = '  +22.05%  '
$ws.Range("E28").Value = '  -4.59%  '
$ws.Range("E29").Value = '  +1.03%  '
$ws.Range("D30").Value = '169.39'
$ws.Range("E30").Value = '  +1.20%  '
$ws.Range("D31").Value = '20.64'
$ws.Range("E31").Value = '  -1.88%  '
$ws.Range("D32").Value = '0.0824'
$ws.Range("E32").Value = '  -4.50%  '
$ws.Range("D33").Value = '0.118'
$ws.Range("E33").Value = '  -6.88%  '
$ws.Range("D34").Value = '29.89'
$ws.Range("E34").Value = '  -4.72%  '
$ws.Range("E35").Value = '  -2.14%  '
$ws.Range("D36").Value = '4.52'
$ws.Range("E36").Value = '  -3.18%  '
$ws.Range("D37").Value = '4.86'
$ws.Range("E37").Value = '  +1.16%  '
$ws.Range("D39").Value = '13.33'
$ws.Range("E39").Value = '  -2.17%  '
$ws.Range("E40").Value = '  -5.54%  '
$ws.Range("D41").Value = '5.77'
$ws.Range("E41").Value = '  -3.30%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '0.202'
$ws.Range("E42").Value = '  -3.94%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").Value = '108.59'
$ws.Range("E43").Value = '  +3.15%  '
$ws.Range("D44").Value = '60.65'
$ws.Range("E44").Value = '  -0.91%  '
$ws.Range("D45").Value = '8.65'
$ws.Range("E45").Value = '  -5.17%  '
$ws.Range("D46").Value = '0.100'
$ws.Range("E46").Value = '  -1.48%  '
$ws.Range("D47").Value = '0.997'
$ws.Range("E47").Value = '  -0.39%  '
$ws.Range("E48").Value = '  -4.08%  '
$ws.Range("E49").Value = '  -1.81%  '
$ws.Range("E50").Value = '  -0.88%  '
$ws.Range("D51").Value = '2.69'
$ws.Range("E51").Value = '  -1.01%  '
